$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Games")
$ws.Rows("2:3").Insert()
Write-Host "done"
